$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in D:E stay text (match source formatting)
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.483.46'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '1.849.54'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('D5').Value = '241.86'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').Value = '0.6267'
$ws.Range('E6').Value = '  -2.65%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('D8').Value = '47.90'
$ws.Range('E8').Value = '  +0.97%  '
$ws.Range('D9').Value = '0.07538'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').Value = '0.2971'
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('D11').Value = '24.26'
$ws.Range('E11').Value = '  -1.12%  '
$ws.Range('D12').Value = '0.07678'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '1.880.03'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '5.010'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '0.6849'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').Value = '83.80'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').Value = '0.000009727'
$ws.Range('E17').Value = '  -2.72%  '
$ws.Range('D18').Value = '2.142.57'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = '6.229'
$ws.Range('E19').Value = '  +1.87%  '
$ws.Range('D20').Value = '29.564.57'
$ws.Range('D21').Value = '234.30'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '7.613'
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('D25').Value = '1.000'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '155.80'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').Value = '0.1388'
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('D28').Value = '8.421'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('D30').Value = '1.482'
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').Value = '0.05833'
$ws.Range('E31').Value = '  -5.67%  '
$ws.Range('D32').Value = '1.260'
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('D33').Value = '4.105'
$ws.Range('E33').Value = '  -1.48%  '
$ws.Range('D34').Value = '4.030'
$ws.Range('E34').Value = '  -1.83%  '
$ws.Range('D35').Value = '1.894'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '1.170'
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').Value = '0.7186'
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('D38').Value = '2.589'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').Value = '2.800'
$ws.Range('E39').Value = '  -0.87%  '
$ws.Range('D40').Value = '1.235.92'
$ws.Range('E40').Value = '  +2.76%  '
$ws.Range('D41').Value = '0.01776'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').Value = '0.9112'
$ws.Range('E42').Value = '  -1.26%  '
$ws.Range('D43').Value = '6.134'
$ws.Range('E43').Value = '  -1.91%  '
$ws.Range('D44').Value = '2.063.21'
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('D45').Value = '0.9999'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '102.85'
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').Value = '67.29'
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('D48').Value = '7.280'
$ws.Range('E48').Value = '  +8.72%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000118'
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.140'
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('E51').Value = '  -0.90%  '

# Restore default (unstyled) cell style now that values are set as text
$priceRange.Style = "Normal"
